$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-11 with revised values (columns B:K)
# Row 2: 2025-03-30
$ws.Range("B2").Value2 = 0.32067577195827568
$ws.Range("C2").Value2 = 0
$ws.Range("D2").Value2 = 0
$ws.Range("E2").Value2 = 0
$ws.Range("F2").Value2 = 0
$ws.Range("G2").Value2 = 0
$ws.Range("H2").Value2 = 0
$ws.Range("I2").Value2 = 0
$ws.Range("J2").Value2 = 0
$ws.Range("K2").Value2 = 0

# Row 3: 2025-04-15
$ws.Range("B3").Value2 = 0.30558067616911933
$ws.Range("C3").Value2 = 0
$ws.Range("D3").Value2 = -0.0024753701868183878
$ws.Range("E3").Value2 = 0.000031947600019851781
$ws.Range("F3").Value2 = -0.00054531602438564911
$ws.Range("G3").Value2 = 0.00054291250795693302
$ws.Range("H3").Value2 = 0.000021548157888519195
$ws.Range("I3").Value2 = -0.0009882259242235058
$ws.Range("J3").Value2 = 0
$ws.Range("K3").Value2 = 0.0029607524759311388

# Row 4: 2025-04-30
$ws.Range("B4").Value2 = 0.29927912833776521
$ws.Range("C4").Value2 = -0.00047674414036869386
$ws.Range("D4").Value2 = 0
$ws.Range("E4").Value2 = 0.000043699683352353059
$ws.Range("F4").Value2 = 0.0000089219695713699257
$ws.Range("G4").Value2 = 0
$ws.Range("H4").Value2 = 0.0000293576986470051
$ws.Range("I4").Value2 = -0.0012389072323672151
$ws.Range("J4").Value2 = -0.000033745392235319395
$ws.Range("K4").Value2 = 0.00010938997100673475

# Row 5: 2025-05-15
$ws.Range("B5").Value2 = 0.27780564978368477
$ws.Range("C5").Value2 = 0.0040149054675811171
$ws.Range("D5").Value2 = -0.0051616445128940276
$ws.Range("E5").Value2 = 0.00014711661762084437
$ws.Range("F5").Value2 = -0.00029367954326683585
$ws.Range("G5").Value2 = -0.0013245540160149033
$ws.Range("H5").Value2 = -0.0001226817332593507
$ws.Range("I5").Value2 = -0.00026105597521251751
$ws.Range("J5").Value2 = 0
$ws.Range("K5").Value2 = -0.0015567056736553453

# Row 6: 2025-05-30
$ws.Range("B6").Value2 = 0.35236365660465274
$ws.Range("C6").Value2 = 0.025148119436511539
$ws.Range("D6").Value2 = 0
$ws.Range("E6").Value2 = -0.00012829680700354224
$ws.Range("F6").Value2 = -0.000087056418044884328
$ws.Range("G6").Value2 = 0
$ws.Range("H6").Value2 = -0.000086576969644260371
$ws.Range("I6").Value2 = -0.0023967068868898616
$ws.Range("J6").Value2 = 0
$ws.Range("K6").Value2 = -0.000038157285548234832

# Row 7: 2025-06-15
$ws.Range("B7").Value2 = 0.34324474238453923
$ws.Range("C7").Value2 = 0
$ws.Range("D7").Value2 = -0.0018358930933242824
$ws.Range("E7").Value2 = -0.00037688357306681005
$ws.Range("F7").Value2 = -0.001767844009046296
$ws.Range("G7").Value2 = 0.00083920263521812825
$ws.Range("H7").Value2 = 0
$ws.Range("I7").Value2 = 0.00015948305304989857
$ws.Range("J7").Value2 = 0
$ws.Range("K7").Value2 = 0.0040001538140037041

# Row 8: 2025-06-30
$ws.Range("B8").Value2 = 0.20301645978423263
$ws.Range("C8").Value2 = -0.034242057974089649
$ws.Range("D8").Value2 = 0
$ws.Range("E8").Value2 = 0.00006456177523307208
$ws.Range("F8").Value2 = -0.00038458769995817279
$ws.Range("G8").Value2 = 0
$ws.Range("H8").Value2 = 0.000038993170911051598
$ws.Range("I8").Value2 = 0.0012944580596345576
$ws.Range("J8").Value2 = 0
$ws.Range("K8").Value2 = 0.00012059182073970165

# Row 9: 2025-07-15
$ws.Range("B9").Value2 = 0.15612360285270505
$ws.Range("C9").Value2 = 0
$ws.Range("D9").Value2 = -0.0013312323546255721
$ws.Range("E9").Value2 = -0.0025542902833823274
$ws.Range("F9").Value2 = -0.0066681728089393513
$ws.Range("G9").Value2 = 0.00077943887493074984
$ws.Range("H9").Value2 = -0.00021125644067366647
$ws.Range("I9").Value2 = -0.00031042344141528274
$ws.Range("J9").Value2 = 0
$ws.Range("K9").Value2 = 0.0003857875314481618

# Row 10: 2025-07-30
$ws.Range("B10").Value2 = 0.38052900968568437
$ws.Range("C10").Value2 = 0.065239114320329666
$ws.Range("D10").Value2 = 0
$ws.Range("E10").Value2 = 0.0000045259273930504009
$ws.Range("F10").Value2 = -0.0004428995244792252
$ws.Range("G10").Value2 = 0
$ws.Range("H10").Value2 = 0.0000048618645659655821
$ws.Range("I10").Value2 = -0.00025564092147157836
$ws.Range("J10").Value2 = -0.0021080902119762213
$ws.Range("K10").Value2 = 0.00013867582425139413

# Row 11: 2025-08-15
$ws.Range("B11").Value2 = 0.42726645926055873
$ws.Range("C11").Value2 = 0
$ws.Range("D11").Value2 = -0.0048355294588895131
$ws.Range("E11").Value2 = 0.0013839191824645137
$ws.Range("F11").Value2 = 0.0022655839343555268
$ws.Range("G11").Value2 = 0.001655503180034068
$ws.Range("H11").Value2 = 0.000055205835676742057
$ws.Range("I11").Value2 = -0.0040860519502526234
$ws.Range("J11").Value2 = 0
$ws.Range("K11").Value2 = 0.026706530002628626

# New row 12: 2025-08-30
$cellA = $ws.Range("A12")
$cellA.NumberFormat = "@"
$cellA.Value2 = "2025-08-30"
$cellA.Style = "Normal"
$ws.Range("B12").Value2 = 0.34111009659532671
$ws.Range("C12").Value2 = -0.040838812509852558
$ws.Range("D12").Value2 = 0
$ws.Range("E12").Value2 = 0.00010144510796152396
$ws.Range("F12").Value2 = -0.0000023485803581274594
$ws.Range("G12").Value2 = 0
$ws.Range("H12").Value2 = -0.00001005758470858199
$ws.Range("I12").Value2 = -0.00099922148469991216
$ws.Range("J12").Value2 = 0
$ws.Range("K12").Value2 = 0.00034681916976642135
